# Fix GLONAF vs GIFT last typo in workflow diagram
#
# This script reproduces, via the PowerPoint COM object model, a set of
# shape repositioning/resizing tweaks plus a text fix on the single
# slide of the workflow diagram.
#
# NOTE on units: Shape.Left/Top/Width/Height are expressed in points,
# while the underlying OOXML <a:off>/<a:ext> are EMUs (1 pt = 12700 EMU).
# The COM layer stores these coordinates as 32-bit floats (Single), so a
# naive "emu / 12700.0" round-trips with +/-1 EMU error after PowerPoint
# truncates back to EMU on save. EmuToPt() nudges the point value by
# tiny increments until the float32 round-trip lands exactly back on the
# requested EMU value, so the saved XML matches byte-for-byte.

function EmuToPt($emu) {
    if ($emu -eq 0) { return 0.0 }
    $sign = 1.0
    $target = $emu
    if ($target -lt 0) {
        $sign = -1.0
        $target = -$target
    }
    $pt = [double]$target / 12700.0
    for ($i = 0; $i -lt 400; $i++) {
        $f = [float]$pt
        $back = [math]::Floor([double]$f * 12700.0)
        if ($back -eq $target) {
            return $sign * $pt
        }
        $pt = $pt + 0.0000002
    }
    return $sign * $pt
}

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $cand = $shapes.Item($i)
        if ($cand.Id -eq $id) {
            return $cand
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shapes = $s.Shapes

# --- Shape id=12 "Rectangle 11" : reposition/resize + text fix ---------
$sh = Get-ShapeById $shapes 12
$sh.Left = EmuToPt(262431)
$sh.Top = EmuToPt(4375586)
$sh.Width = EmuToPt(2682397)
$sh.Height = EmuToPt(1328943)

$tr = $sh.TextFrame.TextRange
$tr.Characters(1, 34).Text = "Alien status for unique species x GIFT"
$tr.Characters(39, 7).Text = " "

# --- Shape id=22 "Rectangle 21" : reposition only -----------------------
$sh = Get-ShapeById $shapes 22
$sh.Left = EmuToPt(8342850)
$sh.Top = EmuToPt(2752047)

# --- Shape id=25 "Straight Arrow Connector 24" : resize (height only) --
$sh = Get-ShapeById $shapes 25
$sh.Height = EmuToPt(699126)

# --- Shape id=27 "Straight Arrow Connector 26" : resize -----------------
$sh = Get-ShapeById $shapes 27
$sh.Width = EmuToPt(2767664)
$sh.Height = EmuToPt(30756)

# --- Shape id=32 "Straight Arrow Connector 31" : resize -----------------
$sh = Get-ShapeById $shapes 32
$sh.Width = EmuToPt(515900)
$sh.Height = EmuToPt(1805494)

# --- Shape id=5 "Rounded Rectangle 4" : reposition -----------------------
$sh = Get-ShapeById $shapes 5
$sh.Left = EmuToPt(289744)
$sh.Top = EmuToPt(2964572)

# --- Shape id=36 "Rounded Rectangle 35" : reposition (y only) -----------
$sh = Get-ShapeById $shapes 36
$sh.Top = EmuToPt(1858204)

# --- Shape id=119 "Straight Arrow Connector 118" : flip + reposition/resize
$sh = Get-ShapeById $shapes 119
$sh.HorizontalFlip = -1
$sh.Left = EmuToPt(9507518)
$sh.Top = EmuToPt(2156014)
$sh.Width = EmuToPt(1)
$sh.Height = EmuToPt(596033)

# --- Shape id=129 "Straight Arrow Connector 128" : reposition/resize ----
$sh = Get-ShapeById $shapes 129
$sh.Left = EmuToPt(1603630)
$sh.Top = EmuToPt(3262382)
$sh.Width = EmuToPt(0)
$sh.Height = EmuToPt(1113204)

# --- Shape id=188 "Connector: Elbow 187" : reposition/resize ------------
$sh = Get-ShapeById $shapes 188
$sh.Left = EmuToPt(3626814)
$sh.Top = EmuToPt(3681344)
$sh.Height = EmuToPt(4760117)
